# This script applies corrected financial data values to rows 2-9
# of the active worksheet, fixing erroneous figures reported previously
# (commit: "error solve ifrs list").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4148
$ws.Range("E2").Value = -286
$ws.Range("F2").Value = -286
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = -19
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = -53
$ws.Range("K2").Value = 5961
$ws.Range("L2").Value = 3034
$ws.Range("M2").Value = 2927
$ws.Range("N2").Value = 2713
$ws.Range("O2").Value = 214
$ws.Range("P2").Value = 196
$ws.Range("Q2").Value = -555
$ws.Range("R2").Value = -195
$ws.Range("S2").Value = 905
$ws.Range("T2").Value = 221
$ws.Range("U2").Value = -776
$ws.Range("V2").Value = 1981
$ws.Range("W2").Value = -6.88
$ws.Range("X2").Value = -0.46
$ws.Range("Y2").Value = 1.25
$ws.Range("Z2").Value = -0.32
$ws.Range("AA2").Value = 103.66
$ws.Range("AB2").Value = 1262.48
$ws.Range("AC2").Value = 83
$ws.Range("AD2").Value = 82.96
$ws.Range("AE2").Value = 6699
$ws.Range("AF2").Value = 1.03
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 40494073

# Row 3
$ws.Range("D3").Value = 3934
$ws.Range("E3").Value = 144
$ws.Range("F3").Value = -308
$ws.Range("G3").Value = 793
$ws.Range("H3").Value = -633
$ws.Range("I3").Value = -383
$ws.Range("J3").Value = -249
$ws.Range("K3").Value = 4367
$ws.Range("L3").Value = 2039
$ws.Range("M3").Value = 2328
$ws.Range("N3").Value = 2331
$ws.Range("O3").Value = -2
$ws.Range("P3").Value = 196
$ws.Range("Q3").Value = 349
$ws.Range("R3").Value = 755
$ws.Range("S3").Value = -1049
$ws.Range("T3").Value = 127
$ws.Range("U3").Value = 221
$ws.Range("V3").Value = 940
$ws.Range("W3").Value = 3.67
$ws.Range("X3").Value = -16.08
$ws.Range("Y3").Value = -15.2
$ws.Range("Z3").Value = -12.25
$ws.Range("AA3").Value = 87.56
$ws.Range("AB3").Value = 1054.46
$ws.Range("AC3").Value = -946
$ws.Range("AD3").Value = -12.28
$ws.Range("AE3").Value = 5755
$ws.Range("AF3").Value = 2.02
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 40494073

# Row 4
$ws.Range("D4").Value = 3990
$ws.Range("E4").Value = 290
$ws.Range("F4").Value = 290
$ws.Range("G4").Value = 203
$ws.Range("H4").Value = 406
$ws.Range("I4").Value = 393
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 3908
$ws.Range("L4").Value = 961
$ws.Range("M4").Value = 2947
$ws.Range("N4").Value = 2864
$ws.Range("O4").Value = 84
$ws.Range("P4").Value = 196
$ws.Range("Q4").Value = 729
$ws.Range("R4").Value = -503
$ws.Range("S4").Value = -235
$ws.Range("T4").Value = 313
$ws.Range("U4").Value = 417
$ws.Range("V4").Value = 132
$ws.Range("W4").Value = 7.27
$ws.Range("X4").Value = 10.18
$ws.Range("Y4").Value = 15.14
$ws.Range("Z4").Value = 9.82
$ws.Range("AA4").Value = 32.6
$ws.Range("AB4").Value = 1255.73
$ws.Range("AC4").Value = 971
$ws.Range("AD4").Value = 13.91
$ws.Range("AE4").Value = 7071
$ws.Range("AF4").Value = 1.91
$ws.Range("AG4").Value = 48
$ws.Range("AH4").Value = 0.36
$ws.Range("AI4").Value = 4.99
$ws.Range("AJ4").Value = 40494073

# Row 5
$ws.Range("D5").Value = 4540
$ws.Range("E5").Value = 495
$ws.Range("F5").Value = 495
$ws.Range("G5").Value = 416
$ws.Range("H5").Value = 422
$ws.Range("I5").Value = 422
$ws.Range("K5").Value = 6503
$ws.Range("L5").Value = 1199
$ws.Range("M5").Value = 5304
$ws.Range("N5").Value = 5304
$ws.Range("P5").Value = 230
$ws.Range("Q5").Value = 612
$ws.Range("R5").Value = -2195
$ws.Range("S5").Value = 1802
$ws.Range("T5").Value = 177
$ws.Range("U5").Value = 435
$ws.Range("V5").Value = 117
$ws.Range("W5").Value = 10.91
$ws.Range("X5").Value = 9.29
$ws.Range("Y5").Value = 10.33
$ws.Range("Z5").Value = 8.11
$ws.Range("AA5").Value = 22.6
$ws.Range("AB5").Value = 2066.05
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 38.3
$ws.Range("AE5").Value = 11531
$ws.Range("AF5").Value = 3.32
$ws.Range("AG5").Value = 650
$ws.Range("AH5").Value = 1.7
$ws.Range("AI5").Value = 70.87
$ws.Range("AJ5").Value = 46000000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 5020
$ws.Range("E6").Value = 487
$ws.Range("F6").Value = 487
$ws.Range("G6").Value = 537
$ws.Range("H6").Value = 413
$ws.Range("I6").Value = 413
$ws.Range("K6").Value = 6922
$ws.Range("L6").Value = 1638
$ws.Range("M6").Value = 5284
$ws.Range("N6").Value = 5284
$ws.Range("P6").Value = 231
$ws.Range("Q6").Value = 860
$ws.Range("R6").Value = -619
$ws.Range("S6").Value = -181
$ws.Range("T6").Value = 1190
$ws.Range("U6").Value = -330
$ws.Range("V6").Value = 391
$ws.Range("W6").Value = 9.71
$ws.Range("X6").Value = 8.23
$ws.Range("Y6").Value = 7.8
$ws.Range("Z6").Value = 6.15
$ws.Range("AA6").Value = 31
$ws.Range("AB6").Value = 2147.87
$ws.Range("AC6").Value = 897
$ws.Range("AD6").Value = 45.71
$ws.Range("AE6").Value = 11460
$ws.Range("AF6").Value = 3.58
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 46110835
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 5783
$ws.Range("E7").Value = 599
$ws.Range("G7").Value = 636
$ws.Range("H7").Value = 491
$ws.Range("I7").Value = 491
$ws.Range("K7").Value = 7997
$ws.Range("L7").Value = 1864
$ws.Range("M7").Value = 6132
$ws.Range("N7").Value = 6132
$ws.Range("P7").Value = 230
$ws.Range("Q7").Value = 708
$ws.Range("R7").Value = -793
$ws.Range("S7").Value = 39
$ws.Range("T7").Value = 1245
$ws.Range("U7").Value = -244
$ws.Range("W7").Value = 10.35
$ws.Range("X7").Value = 8.49
$ws.Range("Y7").Value = 8.6
$ws.Range("Z7").Value = 6.58
$ws.Range("AA7").Value = 30.4
$ws.Range("AC7").Value = 1064
$ws.Range("AD7").Value = 45.11
$ws.Range("AE7").Value = 13298
$ws.Range("AF7").Value = 3.61
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 7439
$ws.Range("E8").Value = 950
$ws.Range("G8").Value = 951
$ws.Range("H8").Value = 724
$ws.Range("I8").Value = 726
$ws.Range("K8").Value = 9345
$ws.Range("L8").Value = 2485
$ws.Range("M8").Value = 6860
$ws.Range("N8").Value = 6860
$ws.Range("P8").Value = 230
$ws.Range("Q8").Value = 1055
$ws.Range("R8").Value = -1123
$ws.Range("S8").Value = 327
$ws.Range("T8").Value = 1213
$ws.Range("U8").Value = 406
$ws.Range("W8").Value = 12.77
$ws.Range("X8").Value = 9.73
$ws.Range("Y8").Value = 11.18
$ws.Range("Z8").Value = 8.35
$ws.Range("AA8").Value = 36.23
$ws.Range("AC8").Value = 1574
$ws.Range("AD8").Value = 30.49
$ws.Range("AE8").Value = 14876
$ws.Range("AF8").Value = 3.23
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 9172
$ws.Range("E9").Value = 1296
$ws.Range("G9").Value = 1278
$ws.Range("H9").Value = 975
$ws.Range("I9").Value = 984
$ws.Range("K9").Value = 10676
$ws.Range("L9").Value = 2837
$ws.Range("M9").Value = 7838
$ws.Range("N9").Value = 7838
$ws.Range("P9").Value = 230
$ws.Range("Q9").Value = 1365
$ws.Range("R9").Value = -964
$ws.Range("S9").Value = 107
$ws.Range("T9").Value = 1000
$ws.Range("U9").Value = 854
$ws.Range("W9").Value = 14.13
$ws.Range("X9").Value = 10.63
$ws.Range("Y9").Value = 13.38
$ws.Range("Z9").Value = 9.74
$ws.Range("AA9").Value = 36.2
$ws.Range("AC9").Value = 2133
$ws.Range("AD9").Value = 22.5
$ws.Range("AE9").Value = 16998
$ws.Range("AF9").Value = 2.82
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
